$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '66.344.53'
$ws.Range('E2').Value2 = '  -1.82%  '

$ws.Range('D3').Value2 = '3.431.82'
$ws.Range('E3').Value2 = '  -4.63%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value2 = '1.00'
$ws.Range('E4').Value2 = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '187.81'
$ws.Range('E5').Value2 = '  -5.75%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '539.72'
$ws.Range('E6').Value2 = '  -3.43%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '0.620'
$ws.Range('E7').Value2 = '  +1.05%  '

$ws.Range('D8').Value2 = '3.424.40'
$ws.Range('E8').Value2 = '  -4.69%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '1.00'
$ws.Range('E9').Value2 = '  -0.08%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '0.640'
$ws.Range('E10').Value2 = '  -4.35%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '58.98'
$ws.Range('E11').Value2 = '  -0.30%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '0.138'
$ws.Range('E12').Value2 = '  -8.97%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '0.0000261'
$ws.Range('E13').Value2 = '  -9.06%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '9.52'
$ws.Range('E14').Value2 = '  -4.60%  '

$ws.Range('D15').Value2 = '3.957.46'
$ws.Range('E15').Value2 = '  -5.44%  '

$ws.Range('E16').Value2 = '  -2.08%  '

$ws.Range('D17').Value2 = '3.407.70'
$ws.Range('E17').Value2 = '  -5.46%  '

$ws.Range('D18').Value2 = '66.044.38'
$ws.Range('E18').Value2 = '  -2.18%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '17.79'
$ws.Range('E19').Value2 = '  -5.99%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '11.49'
$ws.Range('E20').Value2 = '  -6.20%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '1.00'
$ws.Range('E21').Value2 = '  -7.05%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '387.48'
$ws.Range('E22').Value2 = '  -2.78%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '83.94'
$ws.Range('E23').Value2 = '  -1.26%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '3.83'
$ws.Range('E24').Value2 = '  -7.06%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '11.23'
$ws.Range('E25').Value2 = '  -13.05%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '3.91'
$ws.Range('E26').Value2 = '  +1.87%  '

$ws.Range('B27').Value2 = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '11.97'
$ws.Range('E27').Value2 = '  -3.87%  '

$ws.Range('B28').Value2 = 'ImmutableX'
$ws.Range('C28').Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '2.75'
$ws.Range('E28').Value2 = '  -6.56%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '8.73'
$ws.Range('E29').Value2 = '  -7.39%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '689.35'
$ws.Range('E30').Value2 = '  +3.64%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '30.28'
$ws.Range('E31').Value2 = '  -3.67%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '6.95'
$ws.Range('E32').Value2 = '  -18.05%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '11.44'
$ws.Range('E33').Value2 = '  -6.09%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '62.37'
$ws.Range('E34').Value2 = '  -1.97%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '0.108'
$ws.Range('E35').Value2 = '  -4.03%  '

$ws.Range('E36').Value2 = '  +0.10%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '37.49'
$ws.Range('E37').Value2 = '  -11.39%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '0.389'
$ws.Range('E38').Value2 = '  -9.50%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '0.996'
$ws.Range('E39').Value2 = '  -0.26%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '0.128'
$ws.Range('E40').Value2 = '  -5.82%  '

$ws.Range('D41').Value2 = '2.940.27'
$ws.Range('E41').Value2 = '  -9.44%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '2.85'
$ws.Range('E42').Value2 = '  -10.68%  '

$ws.Range('B43').Value2 = 'WEMIXToken'
$ws.Range('C43').Value2 = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '2.73'
$ws.Range('E43').Value2 = '  +0.22%  '

$ws.Range('B44').Value2 = 'PEPE'
$ws.Range('C44').Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D44').Value2 = '0.0₃0643'
$ws.Range('E44').Value2 = '  -16.44%  '

$ws.Range('B45').Value2 = 'Fetch.AI'
$ws.Range('C45').Value2 = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '2.45'
$ws.Range('E45').Value2 = '  -13.12%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '0.0397'
$ws.Range('E46').Value2 = '  -4.57%  '

$ws.Range('E47').Value2 = '  -2.02%  '

$ws.Range('B48').Value2 = 'Monero'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '134.96'
$ws.Range('E48').Value2 = '  -3.00%  '

$ws.Range('B49').Value2 = 'ApeXProtocol'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '2.91'
$ws.Range('E49').Value2 = '  -7.29%  '

$ws.Range('B50').Value2 = 'Stacks'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '2.64'
$ws.Range('E50').Value2 = '  -3.10%  '

$ws.Range('B51').Value2 = 'dogwifhat'
$ws.Range('C51').Value2 = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '2.40'
$ws.Range('E51').Value2 = '  -20.01%  '
